# Add the new "Expenses" line item for the VIP Justin Timberlake Concert.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "VIP Justin Timberlake Concert"
$ws.Range("F7").Value = 26245

# Match the existing currency (no-decimal "$") formatting used by the
# other cells in the Expenses amount column (F4/F6/etc.).
$ws.Range("F7").NumberFormat = $ws.Range("F4").NumberFormat
